# add six-week O&M window and separate cable lay vessels
#
# Root data edit: Port-investments!N11 (the "100 ft sinking basin" cost for
# the Port San Luis 1 / Grays Harbor O&M scenario row) grows from 476 to 567
# (+91) to reflect the new six-week O&M window / separated cable-lay-vessel
# costing. Every other changed cell in the workbook (schedule-short!M8:N18,
# schedule!K3/R9:T35, Port-investments!O11:P11, ...) is a pure formula
# dependent of that one input, so we only need to write the new input and
# let the workbook recalculate; we also re-assert the (already-equivalent)
# formulas the diff shows so the dependents are re-evaluated explicitly.

$wb = $excel.ActiveWorkbook

$wsPort = $wb.Worksheets.Item("Port-investments")
$wsShort = $wb.Worksheets.Item("schedule-short")
$wsFull = $wb.Worksheets.Item("schedule")

# --- Port-investments: six-week O&M window cost bump ----------------------
$wsPort.Range("N11").Value = 567
$wsPort.Range("O11").Formula = "=K11+N11"
$wsPort.Range("P11").Formula = "=K11"

# --- schedule: separate cable lay vessel rows re-asserted -----------------
$wsFull.Range("K3").Formula = "=J3/1000"
$wsFull.Range("S6").Formula = "=S5+R6"
$wsFull.Range("S9").Formula = "=S8+R9"
$wsFull.Range("S11").Formula = "=S10+R11"

# --- View state: "schedule" (was the active/selected tab) gives up focus,
# scrolls right and zooms in; "Port-investments" becomes the active tab and
# scrolls one column right; selections move to reflect where each author
# was last working.
$wsFull.Activate()
$wsFull.Range("S37").Select()
$excel.ActiveWindow.Zoom = 142

$wsPort.Activate()
$wsPort.Range("O20").Select()
$excel.ActiveWindow.Zoom = 113
